$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 row (row 3) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 358
$wsOff.Range("C3").Value = 237
$wsOff.Range("D3").Value = 88
$wsOff.Range("E3").Value = 30
$wsOff.Range("F3").Value = 8

# DEF sheet - Week 17 row (row 3) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 487
$wsDef.Range("C3").Value = 342
$wsDef.Range("D3").Value = 99
$wsDef.Range("E3").Value = 50
